$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the 2021 column (J) header
$ws.Range("J4").Value = 2021

# Add the 2021 data values for column J, rows 5-14
$ws.Range("J5").Value = 24.4
$ws.Range("J6").Value = 45.7
$ws.Range("J7").Value = 38
$ws.Range("J8").Value = 51.3
$ws.Range("J9").Value = 51.5
$ws.Range("J10").Value = 13
$ws.Range("J11").Value = 36.4
$ws.Range("J12").Value = 27
$ws.Range("J13").Value = 2.7
$ws.Range("J14").Value = 40.4

# Copy styles from column I to column J for rows 4-14 so formatting matches
$ws.Range("I4:I14").Copy()
$ws.Range("J4:J14").PasteSpecial(-4122)  # xlPasteFormats

# Adjust row 3 height (18 -> 13.5) per diff
$ws.Rows("3:3").RowHeight = 13.5

# Update the selection to match the target state
$ws.Range("K18").Select()
